$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "27.941.88"
Set-TextValue "E2" "  -1.11%  "

Set-TextValue "D3" "1.749.29"
Set-TextValue "E3" "  -1.94%  "

Set-TextValue "D4" "1.001"

Set-TextValue "D5" "332.51"
Set-TextValue "E5" "  -1.75%  "

Set-TextValue "D6" "0.9978"
Set-TextValue "E6" "  -0.14%  "

Set-TextValue "D7" "0.3879"
Set-TextValue "E7" "  +1.31%  "

Set-TextValue "D8" "0.3381"

Set-TextValue "D9" "45.44"
Set-TextValue "E9" "  -3.44%  "

Set-TextValue "D10" "1.109"
Set-TextValue "E10" "  -3.99%  "

Set-TextValue "D11" "0.07171"
Set-TextValue "E11" "  -3.49%  "

Set-TextValue "D12" "0.9982"

Set-TextValue "D13" "22.09"
Set-TextValue "E13" "  -5.14%  "

Set-TextValue "D14" "6.116"

Set-TextValue "D15" "1.743.10"
Set-TextValue "E15" "  -2.27%  "

Set-TextValue "D16" "6.978"
Set-TextValue "E16" "  -4.43%  "

Set-TextValue "D17" "0.00001049"
Set-TextValue "E17" "  -2.53%  "

Set-TextValue "D18" "0.06596"
Set-TextValue "E18" "  -1.01%  "

Set-TextValue "D19" "80.09"
Set-TextValue "E19" "  -2.95%  "

Set-TextValue "D20" "0.9974"
Set-TextValue "E20" "  -0.22%  "

Set-TextValue "D21" "16.91"
Set-TextValue "E21" "  -3.53%  "

Set-TextValue "D22" "6.186"
Set-TextValue "E22" "  -4.27%  "

Set-TextValue "D23" "27.938.43"
Set-TextValue "E23" "  -1.19%  "

Set-TextValue "E24" "  -4.96%  "

Set-TextValue "D25" "2.386"
Set-TextValue "E25" "  +0.75%  "

Set-TextValue "D26" "153.18"
Set-TextValue "E26" "  -1.56%  "

Set-TextValue "D27" "19.77"
Set-TextValue "E27" "  -5.15%  "

Set-TextValue "D28" "2.294"
Set-TextValue "E28" "  -5.34%  "

Set-TextValue "D29" "1.947.86"
Set-TextValue "E29" "  -1.96%  "

Set-TextValue "D30" "1.260"
Set-TextValue "E30" "  -12.88%  "

Set-TextValue "D31" "128.27"
Set-TextValue "E31" "  -6.78%  "

Set-TextValue "D32" "4.063"
Set-TextValue "E32" "  +3.03%  "

Set-TextValue "D33" "5.781"
Set-TextValue "E33" "  -6.18%  "

Set-TextValue "D34" "0.08609"
Set-TextValue "E34" "  -3.38%  "

Set-TextValue "D35" "11.98"

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D36" "5.097"
Set-TextValue "E36" "  -4.36%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02261"
Set-TextValue "E37" "  -7.35%  "

Set-TextValue "D38" "0.06099"
Set-TextValue "E38" "  -4.19%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D39" "0.6412"
Set-TextValue "E39" "  -6.77%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D40" "1.499"
Set-TextValue "E40" "  +0.16%  "

Set-TextValue "D41" "0.2089"
Set-TextValue "E41" "  -4.24%  "

Set-TextValue "D42" "1.197"
Set-TextValue "E42" "  -3.57%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D43" "0.9968"
Set-TextValue "E43" "  -0.21%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "7.827"
Set-TextValue "E44" "  -5.75%  "

Set-TextValue "D45" "13.65"
Set-TextValue "E45" "  -4.28%  "

Set-TextValue "D46" "3.807"
Set-TextValue "E46" "  -1.66%  "

Set-TextValue "D47" "0.5937"
Set-TextValue "E47" "  -5.90%  "

Set-TextValue "D48" "125.55"
Set-TextValue "E48" "  -5.39%  "

Set-TextValue "D49" "1.967"
Set-TextValue "E49" "  -6.18%  "

Set-TextValue "D50" "0.06979"
Set-TextValue "E50" "  -6.67%  "

Set-TextValue "D51" "1.147"
Set-TextValue "E51" "  -6.48%  "
